# edit.ps1 - applies the RFU3.2 ModificaProfiloAzienda commit:
#  - removes stray w:proofErr spell-check markers
#  - reorders/merges a few runs (email/Partita iva wording tweak)
#  - relocates the _GoBack bookmark to the "numero di telefono" paragraph

$d = $word.ActiveDocument

function Set-ParagraphXml($anchorText, $newParaXml) {
    $r = $d.Content
    $found = $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "anchor not found: $anchorText"
    }
    $pkg = "<?xml version=`"1.0`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" + "<w:body>" + $newParaXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    [void]$r.InsertXML($pkg)
}

Set-ParagraphXml 'RFU 3.2 - ModificaProfiloAzienda' '<w:p w14:paraId="7E5F1043" w14:textId="6F9811BD" w:rsidR="00DB31FD" w:rsidRPr="00DB31FD" w:rsidRDefault="005A0C3C" w:rsidP="00834EAB"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">RFU 3.2 </w:t></w:r><w:r w:rsidR="00BF655F"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r w:rsidR="00FB51BF"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Modifica</w:t></w:r><w:r w:rsidR="0082550F"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>P</w:t></w:r><w:r w:rsidR="00FB51BF"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>rofilo</w:t></w:r><w:r w:rsidR="00953188"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Azienda</w:t></w:r></w:p>'
Set-ParagraphXml 'L’Azienda può modificare il nome' '<w:p w14:paraId="52BBECC0" w14:textId="5F0E266C" w:rsidR="00D27511" w:rsidRPr="00D27511" w:rsidRDefault="00D27511" w:rsidP="00961270"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="368"/></w:pPr><w:r><w:t>L’Azienda può modificare il nome, la via, città, provincia, numero di telefono</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>orario di apertura, orario di chiusura, giorni di apertura. Può modificare la password inserendo la password attuale e la nuova password. Non può modificare Partita iva</w:t></w:r><w:r><w:t xml:space="preserve"> ed email</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>I dati devono essere modificati rispettando i formati specificati nel dizionario dati, sezione DD_Azi</w:t></w:r></w:p>'
Set-ParagraphXml 'ha inserito dati non validi' '<w:p w14:paraId="6A20CDD5" w14:textId="4AEC37C6" w:rsidR="00973C2C" w:rsidRPr="005A7D80" w:rsidRDefault="00961270" w:rsidP="00973C2C"><w:r><w:t>3a</w:t></w:r><w:r w:rsidR="00973C2C" w:rsidRPr="005A7D80"><w:t xml:space="preserve">. Il sistema verifica che </w:t></w:r><w:r w:rsidR="00973C2C"><w:t xml:space="preserve">l’Azienda </w:t></w:r><w:r w:rsidR="00973C2C" w:rsidRPr="005A7D80"><w:t>ha inserito</w:t></w:r><w:r w:rsidR="00973C2C"><w:t xml:space="preserve"> dati non validi, facendo riferimento al dizionario dati, sezione DD_Azi</w:t></w:r></w:p>'
Set-ParagraphXml 'RFU 3.2.a - DatiAziendaModificatiNonValidi' '<w:p w14:paraId="1F30D5FD" w14:textId="15C008E9" w:rsidR="00973C2C" w:rsidRDefault="00973C2C" w:rsidP="00973C2C"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RFU 3.2.a - DatiAziendaModificatiNonValidi</w:t></w:r></w:p>'
Set-ParagraphXml 'Il sistema visualizza un errore di' '<w:p w14:paraId="3903EFF1" w14:textId="0C46B2E9" w:rsidR="00CA607A" w:rsidRPr="005A7D80" w:rsidRDefault="00961270" w:rsidP="00CA607A"><w:r><w:t>4b</w:t></w:r><w:r w:rsidR="00CA607A" w:rsidRPr="005A7D80"><w:t>. Il sistema visualizza un errore di password attuale non corrispondente.</w:t></w:r></w:p>'
Set-ParagraphXml 'PasswordNonCorrispondente' '<w:p w14:paraId="033A8B7D" w14:textId="3B24B892" w:rsidR="00CA607A" w:rsidRPr="00DB31FD" w:rsidRDefault="00CA607A" w:rsidP="00CA607A"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RFU 3.</w:t></w:r><w:r w:rsidR="00973C2C"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> PasswordNonCorrispondente</w:t></w:r></w:p>'

Write-Output "done"
